{"js": "// Find the paragraph that holds the \"When shoud I use WordPad?\" question and\n// append a trailing space to its existing text, then add a new run\n// (identical run formatting) containing the follow-up question, so the\n// paragraph reads:\n//   \"When shoud I use WordPad? Does WordPad have a role to play in\n//    today's office? Did it ever?\"\nconst body = context.document.body;\n\nconst searchText = \"When shoud I use WordPad?\";\nconst results = body.search(searchText, { matchCase: false, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find paragraph text \"${searchText}\" to edit.`);\n}\n\nconst foundRange = results.items[0];\nconst paragraph = foundRange.paragraphs.getFirst();\n\n// 1) Add a trailing space right after the existing question text. Because\n//    this inherits the surrounding run formatting it merges cleanly into the\n//    existing run, turning \"...WordPad?\" into \"...WordPad? \".\nparagraph.insertText(\" \", Word.InsertLocation.end);\nawait context.sync();\n\n// 2) Insert the new sentence as its own run (explicit matching rPr) right\n//    after that, via OOXML so it stays a distinct <w:r> rather than being\n//    silently re-merged with the previous run.\nconst endRange = paragraph.getRange(Word.RangeLocation.end);\nconst newRunOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Calibri\" w:eastAsia=\"Calibri\"/>\n                <w:color w:val=\"auto\"/>\n                <w:spacing w:val=\"0\"/>\n                <w:position w:val=\"0\"/>\n                <w:sz w:val=\"22\"/>\n                <w:shd w:fill=\"auto\" w:val=\"clear\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Does WordPad have a role to play in today's office? Did it ever?</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nendRange.insertOoxml(newRunOoxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Find the paragraph that holds the \"When shoud I use WordPad?\" question and\n# append a trailing space to its existing text, then add a new run\n# (identical run formatting) containing the follow-up question, so the\n# paragraph reads:\n#   \"When shoud I use WordPad? Does WordPad have a role to play in\n#    today's office? Did it ever?\"\n\n$d = $word.ActiveDocument\n\n$searchText = \"When shoud I use WordPad?\"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute($searchText)\n\nif (-not $found) {\n    throw \"Could not find paragraph text '$searchText' to edit.\"\n}\n\n# 1) Add a trailing space right after the existing question text. Because\n#    this inherits the surrounding run formatting it merges cleanly into the\n#    existing run, turning \"...WordPad?\" into \"...WordPad? \".\n$rng.InsertAfter(\" \")\n$rng.Collapse(0)  # wdCollapseEnd\n\n# 2) Insert the new sentence as its own run (explicit matching rPr) right\n#    after that, via OOXML so it stays a distinct <w:r> rather than being\n#    silently re-merged with the previous run.\n$newRunXml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Calibri\" w:eastAsia=\"Calibri\"/>\n                <w:color w:val=\"auto\"/>\n                <w:spacing w:val=\"0\"/>\n                <w:position w:val=\"0\"/>\n                <w:sz w:val=\"22\"/>\n                <w:shd w:fill=\"auto\" w:val=\"clear\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Does WordPad have a role to play in today's office? Did it ever?</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$rng.InsertXML($newRunXml, \"End\")\n"}
